$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027388102864774
$ws.Cells.Item(2, 4).Value = 1.03079118328872
$ws.Cells.Item(2, 5).Value = 1.027386650394401
$ws.Cells.Item(2, 6).Value = 1.025920811995852
$ws.Cells.Item(2, 9).Value = 1.029688338125137
$ws.Cells.Item(2, 10).Value = 1.032546254803232
$ws.Cells.Item(2, 11).Value = 1.033601172550376
$ws.Cells.Item(2, 12).Value = 1.030206529170885
$ws.Cells.Item(2, 13).Value = 1.028744970734982
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029210797770532
$ws.Cells.Item(3, 4).Value = 1.032161511241162
$ws.Cells.Item(3, 5).Value = 1.028964918901997
$ws.Cells.Item(3, 6).Value = 1.028391611944887
$ws.Cells.Item(3, 9).Value = 1.030078471264797
$ws.Cells.Item(3, 10).Value = 1.034004631635953
$ws.Cells.Item(3, 11).Value = 1.034778421433916
$ws.Cells.Item(3, 12).Value = 1.031590430017403
$ws.Cells.Item(3, 13).Value = 1.031018671749184
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030385608663026
$ws.Cells.Item(4, 4).Value = 1.033044086593452
$ws.Cells.Item(4, 5).Value = 1.029982246430839
$ws.Cells.Item(4, 6).Value = 1.029984943621277
$ws.Cells.Item(4, 9).Value = 1.030327797982162
$ws.Cells.Item(4, 10).Value = 1.034943526352196
$ws.Cells.Item(4, 11).Value = 1.035535549548722
$ws.Cells.Item(4, 12).Value = 1.032481525277856
$ws.Cells.Item(4, 13).Value = 1.032484215561182
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030878422396736
$ws.Cells.Item(5, 4).Value = 1.033414151611036
$ws.Cells.Item(5, 5).Value = 1.030409011911067
$ws.Cells.Item(5, 6).Value = 1.030653514811744
$ws.Cells.Item(5, 9).Value = 1.030431875041359
$ws.Cells.Item(5, 10).Value = 1.035337114287145
$ws.Cells.Item(5, 11).Value = 1.035852753030572
$ws.Cells.Item(5, 12).Value = 1.032855110586847
$ws.Cells.Item(5, 13).Value = 1.033099001739449
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030961105517148
$ws.Cells.Item(6, 4).Value = 1.033476230764102
$ws.Cells.Item(6, 5).Value = 1.030480614366127
$ws.Cells.Item(6, 6).Value = 1.030765697720659
$ws.Cells.Item(6, 9).Value = 1.030449306821587
$ws.Cells.Item(6, 10).Value = 1.03540313410845
$ws.Cells.Item(6, 11).Value = 1.035905949230033
$ws.Cells.Item(6, 12).Value = 1.032917777247752
$ws.Cells.Item(6, 13).Value = 1.033202150135097
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030392197863891
$ws.Cells.Item(7, 4).Value = 1.033049035208335
$ws.Cells.Item(7, 5).Value = 1.029987952477596
$ws.Cells.Item(7, 6).Value = 1.029993882024683
$ws.Cells.Item(7, 9).Value = 1.030329191564081
$ws.Cells.Item(7, 10).Value = 1.034948789881744
$ws.Cells.Item(7, 11).Value = 1.035539792311753
$ws.Cells.Item(7, 12).Value = 1.032486521170753
$ws.Cells.Item(7, 13).Value = 1.032492435537804
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028005056527365
$ws.Cells.Item(8, 4).Value = 1.031255154672267
$ws.Cells.Item(8, 5).Value = 1.027920856266567
$ws.Cells.Item(8, 6).Value = 1.02675697953702
$ws.Cells.Item(8, 9).Value = 1.029820834150756
$ws.Cells.Item(8, 10).Value = 1.033040119896044
$ws.Cells.Item(8, 11).Value = 1.033999997460164
$ws.Cells.Item(8, 12).Value = 1.030675142812246
$ws.Cells.Item(8, 13).Value = 1.029514578218143
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.023762369543413
$ws.Cells.Item(9, 4).Value = 1.028061826352144
$ws.Cells.Item(9, 5).Value = 1.024247507504275
$ws.Cells.Item(9, 6).Value = 1.021009788581188
$ws.Cells.Item(9, 9).Value = 1.028900911511854
$ws.Cells.Item(9, 10).Value = 1.029639378597748
$ws.Cells.Item(9, 11).Value = 1.03125051879773
$ws.Cells.Item(9, 12).Value = 1.027448906160816
$ws.Cells.Item(9, 13).Value = 1.024222054496853
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020908066339006
$ws.Cells.Item(10, 4).Value = 1.02591019285647
$ws.Cells.Item(10, 5).Value = 1.021776651489312
$ws.Cells.Item(10, 6).Value = 1.017146806851837
$ws.Cells.Item(10, 9).Value = 1.028271021615735
$ws.Cells.Item(10, 10).Value = 1.02734583275971
$ws.Cells.Item(10, 11).Value = 1.029392232751239
$ws.Cells.Item(10, 12).Value = 1.025273850195779
$ws.Cells.Item(10, 13).Value = 1.020661146287537
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.019665658498755
$ws.Cells.Item(11, 4).Value = 1.024972881107236
$ws.Cells.Item(11, 5).Value = 1.020701263889504
$ws.Cells.Item(11, 6).Value = 1.015466083675844
$ws.Cells.Item(11, 9).Value = 1.027994246508848
$ws.Cells.Item(11, 10).Value = 1.026346168085579
$ws.Cells.Item(11, 11).Value = 1.028581351257984
$ws.Cells.Item(11, 12).Value = 1.024326024664469
$ws.Cells.Item(11, 13).Value = 1.019111017244328
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.019203171090338
$ws.Cells.Item(12, 4).Value = 1.024623853929373
$ws.Cells.Item(12, 5).Value = 1.020300968517793
$ws.Cells.Item(12, 6).Value = 1.014840536117971
$ws.Cells.Item(12, 9).Value = 1.027890827077006
$ws.Cells.Item(12, 10).Value = 1.025973840933575
$ws.Cells.Item(12, 11).Value = 1.028279198215507
$ws.Cells.Item(12, 12).Value = 1.023973034662332
$ws.Cells.Item(12, 13).Value = 1.01853395063955
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.019302422012985
$ws.Cells.Item(13, 4).Value = 1.024698761048725
$ws.Cells.Item(13, 5).Value = 1.020386872015344
$ws.Cells.Item(13, 6).Value = 1.014974775604908
$ws.Cells.Item(13, 9).Value = 1.02791303875128
$ws.Cells.Item(13, 10).Value = 1.026053752338059
$ws.Cells.Item(13, 11).Value = 1.028344054623748
$ws.Cells.Item(13, 12).Value = 1.024048794452074
$ws.Cells.Item(13, 13).Value = 1.018657791988022
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.019627449751082
$ws.Cells.Item(14, 4).Value = 1.024944048235637
$ws.Cells.Item(14, 5).Value = 1.020668192797182
$ws.Cells.Item(14, 6).Value = 1.015414401533084
$ws.Cells.Item(14, 9).Value = 1.027985710372296
$ws.Cells.Item(14, 10).Value = 1.026315412093499
$ws.Cells.Item(14, 11).Value = 1.028556394799374
$ws.Cells.Item(14, 12).Value = 1.024296865408213
$ws.Cells.Item(14, 13).Value = 1.019063343093954
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019827576443929
$ws.Cells.Item(15, 4).Value = 1.025095062148494
$ws.Cells.Item(15, 5).Value = 1.0208414106628
$ws.Cells.Item(15, 6).Value = 1.015685102260009
$ws.Cells.Item(15, 9).Value = 1.02803040432104
$ws.Cells.Item(15, 10).Value = 1.026476495118641
$ws.Cells.Item(15, 11).Value = 1.028687097385053
$ws.Cells.Item(15, 12).Value = 1.024449586817804
$ws.Cells.Item(15, 13).Value = 1.019313045653444
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020990381962374
$ws.Cells.Item(16, 4).Value = 1.025972278576817
$ws.Cells.Item(16, 5).Value = 1.021847903733155
$ws.Cells.Item(16, 6).Value = 1.017258177874437
$ws.Cells.Item(16, 9).Value = 1.028289304745706
$ws.Cells.Item(16, 10).Value = 1.027412037294556
$ws.Cells.Item(16, 11).Value = 1.029445915352658
$ws.Cells.Item(16, 12).Value = 1.025336625705897
$ws.Cells.Item(16, 13).Value = 1.020763846161115
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021718025145026
$ws.Cells.Item(17, 4).Value = 1.026521008864444
$ws.Cells.Item(17, 5).Value = 1.022477763685838
$ws.Cells.Item(17, 6).Value = 1.018242745309833
$ws.Cells.Item(17, 9).Value = 1.028450622256203
$ws.Cells.Item(17, 10).Value = 1.027997109307553
$ws.Cells.Item(17, 11).Value = 1.029920219937332
$ws.Cells.Item(17, 12).Value = 1.025891416660287
$ws.Cells.Item(17, 13).Value = 1.021671660950311
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022141824266523
$ws.Cells.Item(18, 4).Value = 1.026840531156086
$ws.Cells.Item(18, 5).Value = 1.022844622188338
$ws.Cells.Item(18, 6).Value = 1.018816255548458
$ws.Cells.Item(18, 9).Value = 1.028544327845645
$ws.Cells.Item(18, 10).Value = 1.028337742135138
$ws.Cells.Item(18, 11).Value = 1.0301962738641
$ws.Cells.Item(18, 12).Value = 1.026214438288225
$ws.Cells.Item(18, 13).Value = 1.02220038259124
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022286224039145
$ws.Cells.Item(19, 4).Value = 1.026949388546213
$ws.Cells.Item(19, 5).Value = 1.022969622582725
$ws.Cells.Item(19, 6).Value = 1.019011678282565
$ws.Cells.Item(19, 9).Value = 1.028576213414405
$ws.Cells.Item(19, 10).Value = 1.028453783025156
$ws.Cells.Item(19, 11).Value = 1.030290299944441
$ws.Cells.Item(19, 12).Value = 1.026324482849074
$ws.Cells.Item(19, 13).Value = 1.02238053006396
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021640020548739
$ws.Cells.Item(20, 4).Value = 1.02646219158951
$ws.Cells.Item(20, 5).Value = 1.022410240473024
$ws.Cells.Item(20, 6).Value = 1.018137190652185
$ws.Cells.Item(20, 9).Value = 1.028433354620343
$ws.Cells.Item(20, 10).Value = 1.027934401931179
$ws.Cells.Item(20, 11).Value = 1.029869393735215
$ws.Cells.Item(20, 12).Value = 1.025831952842935
$ws.Cells.Item(20, 13).Value = 1.021574343177671
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019531765072958
$ws.Cells.Item(21, 4).Value = 1.024871841396066
$ws.Cells.Item(21, 5).Value = 1.020585374454092
$ws.Cells.Item(21, 6).Value = 1.015284977594884
$ws.Cells.Item(21, 9).Value = 1.027964327369926
$ws.Cells.Item(21, 10).Value = 1.026238387815021
$ws.Cells.Item(21, 11).Value = 1.02849389244507
$ws.Cells.Item(21, 12).Value = 1.024223840420023
$ws.Cells.Item(21, 13).Value = 1.018943954073924
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.01820041029841
$ws.Cells.Item(22, 4).Value = 1.023866893540325
$ws.Cells.Item(22, 5).Value = 1.019433087404771
$ws.Cells.Item(22, 6).Value = 1.013484410090491
$ws.Cells.Item(22, 9).Value = 1.027665881804098
$ws.Cells.Item(22, 10).Value = 1.02516619700473
$ws.Cells.Item(22, 11).Value = 1.02762352183618
$ws.Cells.Item(22, 12).Value = 1.023207390995284
$ws.Cells.Item(22, 13).Value = 1.017282698574559
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.018906747236364
$ws.Cells.Item(23, 4).Value = 1.024400119274367
$ws.Cells.Item(23, 5).Value = 1.020044411050158
$ws.Cells.Item(23, 6).Value = 1.014439629552488
$ws.Cells.Item(23, 9).Value = 1.027824432463408
$ws.Cells.Item(23, 10).Value = 1.025735147264079
$ws.Cells.Item(23, 11).Value = 1.02808545328471
$ws.Cells.Item(23, 12).Value = 1.023746746136102
$ws.Cells.Item(23, 13).Value = 1.018164079893409
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021675269391704
$ws.Cells.Item(24, 4).Value = 1.026488770260924
$ws.Cells.Item(24, 5).Value = 1.022440752937626
$ws.Cells.Item(24, 6).Value = 1.018184888639803
$ws.Cells.Item(24, 9).Value = 1.028441158320963
$ws.Cells.Item(24, 10).Value = 1.027962738640553
$ws.Cells.Item(24, 11).Value = 1.029892361761262
$ws.Cells.Item(24, 12).Value = 1.025858823769387
$ws.Cells.Item(24, 13).Value = 1.021618319333701
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.024863651307452
$ws.Cells.Item(25, 4).Value = 1.028891310207258
$ws.Cells.Item(25, 5).Value = 1.025200936565805
$ws.Cells.Item(25, 6).Value = 1.022500954845393
$ws.Cells.Item(25, 9).Value = 1.029141632367965
$ws.Cells.Item(25, 10).Value = 1.030523111253799
$ws.Cells.Item(25, 11).Value = 1.031965711012572
$ws.Cells.Item(25, 12).Value = 1.028287152545051
$ws.Cells.Item(25, 13).Value = 1.025595874297711
